$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update res_bus vm_pu values for the 380 kV case (Case_2_9)
$data = @{
    "B2"=1.02; "C2"=1.024813526207424; "D2"=1.034954486447669; "E2"=1.045577046498844; "F2"=1.048882581683446; "I2"=1.032431564161194; "J2"=1.029986039473223; "K2"=1.037752478194803; "L2"=1.048344896668772; "M2"=1.051641187459325; "N2"=1.014048805542006;
    "B3"=1.02; "C3"=1.02564705675749; "D3"=1.035599000728249; "E3"=1.046461160054676; "F3"=1.049738042079583; "I3"=1.032561389820101; "J3"=1.030459236348806; "K3"=1.038206725891144; "L3"=1.049040289677001; "M3"=1.052308670455548; "N3"=1.014205319681051;
    "B4"=1.02; "C4"=1.02618672696388; "D4"=1.036015951258566; "E4"=1.04703433821752; "F4"=1.050292253851797; "I4"=1.032643570670792; "J4"=1.030765101555332; "K4"=1.038499875618053; "L4"=1.049490714784494; "M4"=1.052740604221164; "N4"=1.014306471700703;
    "B5"=1.02; "C5"=1.026413679451179; "D5"=1.036191213091241; "E5"=1.047275563276359; "F5"=1.050525403403729; "I5"=1.032677681673764; "J5"=1.030893608301211; "K5"=1.038622927877364; "L5"=1.049680181925941; "M5"=1.052922194170008; "N5"=1.014348966111216;
    "B6"=1.02; "C6"=1.026451790133882; "D6"=1.036220638853349; "E6"=1.04731608127551; "F6"=1.050564559503928; "I6"=1.032683383360328; "J6"=1.03091518048175; "K6"=1.038643577834213; "L6"=1.049712000606064; "M6"=1.052952684167054; "N6"=1.014356099343942;
    "B7"=1.02; "C7"=1.026189759220818; "D7"=1.036018293213043; "E7"=1.04703756045544; "F7"=1.050295368585435; "I7"=1.032644028185193; "J7"=1.030766818980671; "K7"=1.038501520588415; "L7"=1.049493246029609; "M7"=1.052743030617371; "N7"=1.014307039631119;
    "B8"=1.02; "C8"=1.025095154542803; "D8"=1.035172321286735; "E8"=1.045875608471968; "F8"=1.049171548321055; "I8"=1.032475816689571; "J8"=1.030146025091863; "K8"=1.037906153760223; "L8"=1.048579811959865; "M8"=1.051866759767541; "N8"=1.014101725457729;
    "B9"=1.02; "C9"=1.02316884084701; "D9"=1.033680976691155; "E9"=1.04383658624417; "F9"=1.047196460750341; "I9"=1.032165468939527; "J9"=1.029049673219147; "K9"=1.036851138025292; "L9"=1.046973810215493; "M9"=1.050322941670483; "N9"=1.013739011807179;
    "B10"=1.02; "C10"=1.021886415571814; "D10"=1.032686428945423; "E10"=1.042483050172753; "F10"=1.045883361586555; "I10"=1.031949252993863; "J10"=1.028317207943732; "K10"=1.036143908444685; "L10"=1.045905644370269; "M10"=1.049294015785099; "N10"=1.013496606909629;
    "B11"=1.02; "C11"=1.021331553063862; "D11"=1.032255726048248; "E11"=1.041898355078425; "F11"=1.045315658492637; "I11"=1.031853430235042; "J11"=1.027999686519722; "K11"=1.03583676597853; "L11"=1.045443729301334; "M11"=1.048848567456917; "N11"=1.013391506896466;
    "B12"=1.02; "C12"=1.021125519419062; "D12"=1.032095737042351; "E12"=1.041681384667231; "F12"=1.045104921992129; "I12"=1.031817507757005; "J12"=1.027881692167234; "K12"=1.035722544643271; "L12"=1.045272246321278; "M12"=1.048683122239674; "N12"=1.013352447896208;
    "B13"=1.02; "C13"=1.021169711279862; "D13"=1.032130055504638; "E13"=1.041727915962611; "F13"=1.045150119606416; "I13"=1.031825228153575; "J13"=1.027907004724619; "K13"=1.035747051592034; "L13"=1.045309025776726; "M13"=1.048718610135537; "N13"=1.013360827090607;
    "B14"=1.02; "C14"=1.021314520883276; "D14"=1.032242501434807; "E14"=1.041880415914032; "F14"=1.045298236209414; "I14"=1.031850467591376; "J14"=1.027989934139076; "K14"=1.035827327152275; "L14"=1.045429552561896; "M14"=1.048834891406255; "N14"=1.013388278676471;
    "B15"=1.02; "C15"=1.02140375178103; "D15"=1.032311782240255; "E15"=1.041974404231104; "F15"=1.045389513510126; "I15"=1.031865974787315; "J15"=1.028041022733128; "K15"=1.035876769733575; "L15"=1.045503825440472; "M15"=1.048906538054458; "N15"=1.013405189843615;
    "B16"=1.02; "C16"=1.021923249283271; "D16"=1.03271501222901; "E16"=1.042521884004689; "F16"=1.045921056819015; "I16"=1.031955566160479; "J16"=1.028338273330768; "K16"=1.036164273478636; "L16"=1.045936313088196; "M16"=1.049323580616772; "N16"=1.013503579198932;
    "B17"=1.02; "C17"=1.022249234140602; "D17"=1.03296793372794; "E17"=1.042865678346218; "F17"=1.046254716148861; "I17"=1.032011176114363; "J17"=1.028524635475119; "K17"=1.036344375022071; "L17"=1.046207764966995; "M17"=1.049585203849358; "N17"=1.013565259836014;
    "B18"=1.02; "C18"=1.022439417570909; "D18"=1.033115452922758; "E18"=1.043066342241344; "F18"=1.046449418532646; "I18"=1.032043400150718; "J18"=1.028633302630409; "K18"=1.036449337563179; "L18"=1.046366156795873; "M18"=1.049737812191542; "N18"=1.013601223852451;
    "B19"=1.02; "C19"=1.022504272281149; "D19"=1.033165752144592; "E19"=1.043134786197868; "F19"=1.046515821296759; "I19"=1.032054351676857; "J19"=1.028670349394792; "K19"=1.036485112169612; "L19"=1.046420174199231; "M19"=1.049789848989976; "N19"=1.013613484391454;
    "B20"=1.02; "C20"=1.022214254697803; "D20"=1.032940798214603; "E20"=1.042828778544876; "F20"=1.046218908902419; "I20"=1.032005231646471; "J20"=1.028504644164423; "K20"=1.036325060886188; "L20"=1.046178634683209; "M20"=1.049557133319952; "N20"=1.013558643455441;
    "B21"=1.02; "C21"=1.021271876203085; "D21"=1.03220938908344; "E21"=1.041835502629692; "F21"=1.045254615855618; "I21"=1.03184304429878; "J21"=1.027965514940266; "K21"=1.035803691714786; "L21"=1.045394057852514; "M21"=1.048800649067503; "N21"=1.013380195422225;
    "B22"=1.02; "C22"=1.020679753788637; "D22"=1.031749484840475; "E22"=1.041212214511702; "F22"=1.044649101046103; "I22"=1.031739163909749; "J22"=1.027626238347454; "K22"=1.035475106135013; "L22"=1.044901301177362; "M22"=1.048325099907742; "N22"=1.01326788153757;
    "B23"=1.02; "C23"=1.020993611780043; "D23"=1.031993291807043; "E23"=1.041542514725349; "F23"=1.04497002186629; "I23"=1.031794413322768; "J23"=1.027806123767463; "K23"=1.035649369052835; "L23"=1.045162469315647; "M23"=1.048577189160387; "N23"=1.013327432157805;
    "B24"=1.02; "C24"=1.022230060272775; "D24"=1.032953059602669; "E24"=1.04284545156148; "F24"=1.046235088394773; "I24"=1.032007918351718; "J24"=1.028513677484651; "K24"=1.036333788383774; "L24"=1.046191797223845; "M24"=1.049569817159562; "N24"=1.013561633154027;
    "B25"=1.02; "C25"=1.023666531263994; "D25"=1.034066588963819; "E25"=1.044362705757357; "F25"=1.047706437277641; "I25"=1.032247347427805; "J25"=1.029333387404381; "K25"=1.037124575848873; "L25"=1.047388565829344; "M25"=1.050722012213266; "N25"=1.013832888843287
}

foreach ($addr in $data.Keys) {
    $ws.Range($addr).Value = $data[$addr]
}

